# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect the latest scrape, per commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# 展览 sheet: rows 5, 8, 11, 35
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 15548
$wsExpo.Range("F8").Value = 702
$wsExpo.Range("F11").Value = 8991
$wsExpo.Range("F35").Value = 316

# 全部类型 sheet: rows 5, 8, 11, 37
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 15548
$wsAll.Range("F8").Value = 702
$wsAll.Range("F11").Value = 8991
$wsAll.Range("F37").Value = 316
